# Auto-generated: apply scheduled-runner price/profit updates to Kujata_Profits workbook
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 1553.35
$ws.Cells.Item(70, 10).Value = 1652.75
$ws.Cells.Item(70, 12).Value = 4958.25
$ws.Cells.Item(70, 14).Value = -5498.25
$ws.Cells.Item(73, 8).Value = 1553.35
$ws.Cells.Item(73, 10).Value = 1652.75
$ws.Cells.Item(73, 12).Value = 4958.25
$ws.Cells.Item(73, 14).Value = -6830.25
$ws.Cells.Item(80, 8).Value = 798.875
$ws.Cells.Item(80, 9).Value = 1258.8
$ws.Cells.Item(80, 11).Value = 3776.4
$ws.Cells.Item(80, 13).Value = -2778.4
$ws.Cells.Item(83, 8).Value = 798.875
$ws.Cells.Item(83, 9).Value = 1258.8
$ws.Cells.Item(83, 11).Value = 11329.2
$ws.Cells.Item(83, 13).Value = -6337.199999999999
$ws.Cells.Item(88, 8).Value = 823641.8
$ws.Cells.Item(88, 9).Value = 314
$ws.Cells.Item(88, 10).Value = 1123033.8
$ws.Cells.Item(88, 11).Value = 314
$ws.Cells.Item(88, 12).Value = 1123033.8
$ws.Cells.Item(88, 13).Value = 92
$ws.Cells.Item(88, 14).Value = -1123845.8
$ws.Cells.Item(91, 8).Value = 823641.8
$ws.Cells.Item(91, 9).Value = 314
$ws.Cells.Item(91, 10).Value = 1123033.8
$ws.Cells.Item(91, 11).Value = 314
$ws.Cells.Item(91, 12).Value = 1123033.8
$ws.Cells.Item(91, 13).Value = 1090
$ws.Cells.Item(91, 14).Value = -1125841.8
$ws.Cells.Item(98, 8).Value = 4628.923
$ws.Cells.Item(98, 9).Value = 3039.3914
$ws.Cells.Item(98, 10).Value = 16815.334
$ws.Cells.Item(98, 11).Value = 3039.3914
$ws.Cells.Item(98, 12).Value = 16815.334
$ws.Cells.Item(98, 13).Value = -1541.3914
$ws.Cells.Item(98, 14).Value = -19811.334
$ws.Cells.Item(100, 8).Value = 1692.1428
$ws.Cells.Item(100, 9).Value = 1692.1428
$ws.Cells.Item(100, 11).Value = 1692.1428
$ws.Cells.Item(100, 13).Value = -1151.1428
$ws.Cells.Item(112, 8).Value = 2336.9714
$ws.Cells.Item(112, 10).Value = 2648.0688
$ws.Cells.Item(112, 12).Value = 7944.2064
$ws.Cells.Item(112, 14).Value = -10160.2064
$ws.Cells.Item(122, 8).Value = 4628.923
$ws.Cells.Item(122, 9).Value = 3039.3914
$ws.Cells.Item(122, 10).Value = 16815.334
$ws.Cells.Item(122, 11).Value = 9118.174199999999
$ws.Cells.Item(122, 12).Value = 50446.00199999999
$ws.Cells.Item(122, 13).Value = -6668.174199999999
$ws.Cells.Item(122, 14).Value = -55346.00199999999
$ws.Cells.Item(138, 8).Value = 1509.89
$ws.Cells.Item(138, 9).Value = 780.25714
$ws.Cells.Item(138, 10).Value = 1902.7693
$ws.Cells.Item(138, 11).Value = 2340.77142
$ws.Cells.Item(138, 12).Value = 5708.3079
$ws.Cells.Item(138, 13).Value = 2799.22858
$ws.Cells.Item(138, 14).Value = -15988.3079

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 1595.7778
$ws.Cells.Item(61, 9).Value = 1272.4
$ws.Cells.Item(61, 11).Value = 1272.4
$ws.Cells.Item(61, 13).Value = -1060.4
$ws.Cells.Item(88, 8).Value = 2612
$ws.Cells.Item(88, 9).Value = 2299.5
$ws.Cells.Item(88, 10).Value = 2660.077
$ws.Cells.Item(88, 11).Value = 2299.5
$ws.Cells.Item(88, 12).Value = 2660.077
$ws.Cells.Item(88, 13).Value = -1893.5
$ws.Cells.Item(88, 14).Value = -3472.077
$ws.Cells.Item(91, 8).Value = 2612
$ws.Cells.Item(91, 9).Value = 2299.5
$ws.Cells.Item(91, 10).Value = 2660.077
$ws.Cells.Item(91, 11).Value = 2299.5
$ws.Cells.Item(91, 12).Value = 2660.077
$ws.Cells.Item(91, 13).Value = -895.5
$ws.Cells.Item(91, 14).Value = -5468.077
$ws.Cells.Item(136, 8).Value = 1595.7778
$ws.Cells.Item(136, 9).Value = 1272.4
$ws.Cells.Item(136, 11).Value = 3817.2
$ws.Cells.Item(136, 13).Value = -1267.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3957.6924
$ws.Cells.Item(86, 9).Value = 4158
$ws.Cells.Item(86, 10).Value = 3414
$ws.Cells.Item(86, 11).Value = 4158
$ws.Cells.Item(86, 12).Value = 3414
$ws.Cells.Item(86, 13).Value = -3035
$ws.Cells.Item(86, 14).Value = -5660
$ws.Cells.Item(89, 8).Value = 3957.6924
$ws.Cells.Item(89, 9).Value = 4158
$ws.Cells.Item(89, 10).Value = 3414
$ws.Cells.Item(89, 11).Value = 20790
$ws.Cells.Item(89, 12).Value = 17070
$ws.Cells.Item(89, 13).Value = -15174
$ws.Cells.Item(89, 14).Value = -28302
$ws.Cells.Item(102, 8).Value = 33166.5
$ws.Cells.Item(102, 9).Value = 24000
$ws.Cells.Item(102, 11).Value = 24000
$ws.Cells.Item(102, 13).Value = -20755
$ws.Cells.Item(134, 8).Value = 8130.2354
$ws.Cells.Item(134, 9).Value = 1185.4166
$ws.Cells.Item(134, 10).Value = 24797.8
$ws.Cells.Item(134, 11).Value = 3556.2498
$ws.Cells.Item(134, 12).Value = 74393.39999999999
$ws.Cells.Item(134, 13).Value = -1021.2498
$ws.Cells.Item(134, 14).Value = -79463.39999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1453.5862
$ws.Cells.Item(31, 9).Value = 1317.804
$ws.Cells.Item(31, 11).Value = 1317.804
$ws.Cells.Item(31, 13).Value = -1022.804
$ws.Cells.Item(34, 8).Value = 1453.5862
$ws.Cells.Item(34, 9).Value = 1317.804
$ws.Cells.Item(34, 11).Value = 1317.804
$ws.Cells.Item(34, 13).Value = -1115.804
$ws.Cells.Item(132, 8).Value = 2413
$ws.Cells.Item(132, 9).Value = 1784.7858
$ws.Cells.Item(132, 11).Value = 5354.357400000001
$ws.Cells.Item(132, 13).Value = -2824.357400000001
$ws.Cells.Item(134, 8).Value = 1361.174
$ws.Cells.Item(134, 9).Value = 1228.9048
$ws.Cells.Item(134, 11).Value = 3686.7144
$ws.Cells.Item(134, 13).Value = -1151.7144

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1063652.5
$ws.Cells.Item(4, 9).Value = 100019.78
$ws.Cells.Item(4, 11).Value = 300059.34
$ws.Cells.Item(4, 13).Value = -299947.34
$ws.Cells.Item(107, 8).Value = 17266.5
$ws.Cells.Item(107, 10).Value = 20519.8
$ws.Cells.Item(107, 12).Value = 61559.39999999999
$ws.Cells.Item(107, 14).Value = -65399.39999999999
$ws.Cells.Item(116, 8).Value = 2666.6667
$ws.Cells.Item(116, 9).Value = 1000
$ws.Cells.Item(116, 11).Value = 3000
$ws.Cells.Item(116, 13).Value = 442
$ws.Cells.Item(131, 8).Value = 13160125
$ws.Cells.Item(131, 10).Value = 2628.0952
$ws.Cells.Item(131, 12).Value = 7884.285600000001
$ws.Cells.Item(131, 14).Value = -17964.2856
$ws.Cells.Item(136, 8).Value = 2160.3333
$ws.Cells.Item(136, 10).Value = 1666.4
$ws.Cells.Item(136, 12).Value = 4999.200000000001
$ws.Cells.Item(136, 14).Value = -15199.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 6796.6665
$ws.Cells.Item(80, 9).Value = 6790
$ws.Cells.Item(80, 11).Value = 6790
$ws.Cells.Item(80, 13).Value = -5792
$ws.Cells.Item(83, 8).Value = 6796.6665
$ws.Cells.Item(83, 9).Value = 6790
$ws.Cells.Item(83, 11).Value = 33950
$ws.Cells.Item(83, 13).Value = -28958
$ws.Cells.Item(102, 8).Value = 4856
$ws.Cells.Item(102, 9).Value = 3798.4443
$ws.Cells.Item(102, 10).Value = 5450.875
$ws.Cells.Item(102, 11).Value = 3798.4443
$ws.Cells.Item(102, 12).Value = 5450.875
$ws.Cells.Item(102, 13).Value = -2176.4443
$ws.Cells.Item(102, 14).Value = -8694.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1700.2727
$ws.Cells.Item(68, 9).Value = 1475.625
$ws.Cells.Item(68, 11).Value = 1475.625
$ws.Cells.Item(68, 13).Value = -726.625
$ws.Cells.Item(71, 8).Value = 1700.2727
$ws.Cells.Item(71, 9).Value = 1475.625
$ws.Cells.Item(71, 11).Value = 7378.125
$ws.Cells.Item(71, 13).Value = -3634.125
$ws.Cells.Item(93, 8).Value = 901.5925999999999
$ws.Cells.Item(93, 9).Value = 623.6087
$ws.Cells.Item(93, 11).Value = 623.6087
$ws.Cells.Item(93, 13).Value = 624.3913
$ws.Cells.Item(132, 8).Value = 32106.424
$ws.Cells.Item(132, 9).Value = 1119
$ws.Cells.Item(132, 10).Value = 61271.06
$ws.Cells.Item(132, 11).Value = 3357
$ws.Cells.Item(132, 12).Value = 183813.18
$ws.Cells.Item(132, 13).Value = -827
$ws.Cells.Item(132, 14).Value = -188873.18

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(14, 8).Value = 234001.67
$ws.Cells.Item(14, 9).Value = 350502.5
$ws.Cells.Item(14, 11).Value = 350502.5
$ws.Cells.Item(14, 13).Value = -350334.5
$ws.Cells.Item(81, 8).Value = 420.8
$ws.Cells.Item(81, 9).Value = 420.8
$ws.Cells.Item(81, 10).Value = 0
$ws.Cells.Item(81, 11).Value = 841.6
$ws.Cells.Item(81, 12).Value = 0
$ws.Cells.Item(81, 13).Value = 219.4
$ws.Cells.Item(81, 14).ClearContents()
$ws.Cells.Item(84, 8).Value = 420.8
$ws.Cells.Item(84, 9).Value = 420.8
$ws.Cells.Item(84, 10).Value = 0
$ws.Cells.Item(84, 11).Value = 4208
$ws.Cells.Item(84, 12).Value = 0
$ws.Cells.Item(84, 13).Value = 1096
$ws.Cells.Item(84, 14).ClearContents()
$ws.Cells.Item(136, 8).Value = 684.8570999999999
$ws.Cells.Item(136, 9).Value = 632.3333
$ws.Cells.Item(136, 11).Value = 1896.9999
$ws.Cells.Item(136, 13).Value = 653.0001
